# Regenerate merged AHB files
# - rename the "_old" / "_new" header suffixes to "_FV2304" / "_FV2310"
# - wrap the data range in a native Excel Table (ListObject)
# - freeze the header row (split at row 2, frozen)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$baseNames = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")

# Columns A-J (1-10): "<Name>_old" -> "<Name>_FV2304"
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $baseNames[$i] + "_FV2304"
}

# Column K (11) stays "diff"

# Columns L-U (12-21): "<Name>_new" -> "<Name>_FV2310"
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $baseNames[$i] + "_FV2310"
}

# Turn the used range into a proper Excel Table, header row included.
$dataRange = $ws.Range("A1:U82")
$tbl = $ws.ListObjects.Add(1, $dataRange, 0, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# Freeze panes above row 2 (i.e. freeze the header row).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select()
